$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1494.3
$ws.Range("I6").Value = 1290
$ws.Range("J6").Value = 3333
$ws.Range("K6").Value = 3870
$ws.Range("L6").Value = 9999
$ws.Range("M6").Value = -3758
$ws.Range("N6").Value = -10223
# Row 8
$ws.Range("H8").Value = 55.416668
$ws.Range("I8").Value = 55.416668
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 166.250004
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -27.25000399999999
$ws.Range("N8").ClearContents()
# Row 70
$ws.Range("H70").Value = 1092.6
$ws.Range("I70").Value = 1063.5
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 3190.5
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -2920.5
$ws.Range("N70").Value = -5040
# Row 73
$ws.Range("H73").Value = 1092.6
$ws.Range("I73").Value = 1063.5
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 3190.5
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2254.5
$ws.Range("N73").Value = -6372
# Row 76
$ws.Range("H76").Value = 3071.3333
$ws.Range("I76").Value = 3059.6365
$ws.Range("K76").Value = 3059.6365
$ws.Range("M76").Value = -2744.6365
# Row 79
$ws.Range("H79").Value = 3071.3333
$ws.Range("I79").Value = 3059.6365
$ws.Range("K79").Value = 3059.6365
$ws.Range("M79").Value = -1967.6365
# Row 106
$ws.Range("H106").Value = 112892.43
$ws.Range("I106").Value = 1992.5
$ws.Range("J106").Value = 260759
$ws.Range("K106").Value = 1992.5
$ws.Range("L106").Value = 260759
$ws.Range("M106").Value = -1361.5
$ws.Range("N106").Value = -262021
# Row 138
$ws.Range("H138").Value = 2117.1836
$ws.Range("I138").Value = 1684.4166
$ws.Range("J138").Value = 2368.4678
$ws.Range("K138").Value = 5053.2498
$ws.Range("L138").Value = 7105.403399999999
$ws.Range("M138").Value = 86.7502000000004
$ws.Range("N138").Value = -17385.4034

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 27311.086
$ws.Range("I32").Value = 26310.09
$ws.Range("J32").Value = 49666.668
$ws.Range("K32").Value = 26310.09
$ws.Range("L32").Value = 49666.668
$ws.Range("M32").Value = -26023.09
$ws.Range("N32").Value = -50240.668
# Row 63
$ws.Range("H63").Value = 2928.2354
$ws.Range("I63").Value = 2473.3333
$ws.Range("J63").Value = 4020
$ws.Range("K63").Value = 2473.3333
$ws.Range("L63").Value = 4020
$ws.Range("M63").Value = -1787.3333
$ws.Range("N63").Value = -5392
# Row 66
$ws.Range("H66").Value = 2928.2354
$ws.Range("I66").Value = 2473.3333
$ws.Range("J66").Value = 4020
$ws.Range("K66").Value = 12366.6665
$ws.Range("L66").Value = 20100
$ws.Range("M66").Value = -8934.666499999999
$ws.Range("N66").Value = -26964
# Row 80
$ws.Range("H80").Value = 54672.8
$ws.Range("J80").Value = 54672.8
$ws.Range("L80").Value = 54672.8
$ws.Range("N80").Value = -56668.8
# Row 83
$ws.Range("H83").Value = 54672.8
$ws.Range("J83").Value = 54672.8
$ws.Range("L83").Value = 164018.4
$ws.Range("N83").Value = -174002.4
# Row 120
$ws.Range("H120").Value = 45058
$ws.Range("J120").Value = 45058
$ws.Range("L120").Value = 45058
$ws.Range("N120").Value = -54734
# Row 122
$ws.Range("H122").Value = 3183.8333
$ws.Range("I122").Value = 4407.6665
$ws.Range("K122").Value = 13222.9995
$ws.Range("M122").Value = -10772.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 112
$ws.Range("H112").Value = 46900
$ws.Range("J112").Value = 46900
$ws.Range("L112").Value = 46900
$ws.Range("N112").Value = -49854
# Row 130
$ws.Range("H130").Value = 45492.285
$ws.Range("J130").Value = 45492.285
$ws.Range("L130").Value = 45492.285
$ws.Range("N130").Value = -55532.285

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2137.392
$ws.Range("I58").Value = 1857.119
$ws.Range("J58").Value = 3445.3333
$ws.Range("K58").Value = 1857.119
$ws.Range("L58").Value = 3445.3333
$ws.Range("M58").Value = -1654.119
$ws.Range("N58").Value = -3851.3333
# Row 110
$ws.Range("H110").Value = 41350.5
$ws.Range("J110").Value = 41350.5
$ws.Range("L110").Value = 41350.5
$ws.Range("N110").Value = -49530.5
# Row 132
$ws.Range("H132").Value = 31646.596
$ws.Range("I132").Value = 1378.7894
$ws.Range("K132").Value = 4136.3682
$ws.Range("M132").Value = -1606.3682
# Row 134
$ws.Range("H134").Value = 3063.7222
$ws.Range("I134").Value = 1472
$ws.Range("J134").Value = 4655.4443
$ws.Range("K134").Value = 4416
$ws.Range("L134").Value = 13966.3329
$ws.Range("M134").Value = -1881
$ws.Range("N134").Value = -19036.3329
# Row 136
$ws.Range("H136").Value = 2137.392
$ws.Range("I136").Value = 1857.119
$ws.Range("J136").Value = 3445.3333
$ws.Range("K136").Value = 5571.357
$ws.Range("L136").Value = 10335.9999
$ws.Range("M136").Value = -3021.357
$ws.Range("N136").Value = -15435.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 9709.727999999999
$ws.Range("I107").Value = 10583.6
$ws.Range("J107").Value = 8981.5
$ws.Range("K107").Value = 31750.8
$ws.Range("L107").Value = 26944.5
$ws.Range("M107").Value = -29830.8
$ws.Range("N107").Value = -30784.5
# Row 112
$ws.Range("H112").Value = 2988410.2
$ws.Range("I112").Value = 66668000
$ws.Range("J112").Value = 3429.5312
$ws.Range("K112").Value = 200004000
$ws.Range("L112").Value = 10288.5936
$ws.Range("M112").Value = -200002892
$ws.Range("N112").Value = -12504.5936
# Row 122
$ws.Range("H122").Value = 2900.6382
$ws.Range("I122").Value = 688.3171
$ws.Range("J122").Value = 18018.166
$ws.Range("K122").Value = 6194.8539
$ws.Range("L122").Value = 162163.494
$ws.Range("M122").Value = -3744.8539
$ws.Range("N122").Value = -167063.494
# Row 137
$ws.Range("H137").Value = 10069
$ws.Range("I137").Value = 3200
$ws.Range("J137").Value = 20372.5
$ws.Range("K137").Value = 9600
$ws.Range("L137").Value = 61117.5
$ws.Range("M137").Value = -4500
$ws.Range("N137").Value = -71317.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 24199
$ws.Range("I26").Value = 24199
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 24199
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -23919
$ws.Range("N26").ClearContents()
# Row 50
$ws.Range("H50").Value = 24199
$ws.Range("I50").Value = 24199
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 24199
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -23701
$ws.Range("N50").ClearContents()
# Row 58
$ws.Range("H58").Value = 22000
$ws.Range("J58").Value = 22000
$ws.Range("L58").Value = 22000
$ws.Range("N58").Value = -22554
# Row 97
$ws.Range("H97").Value = 4915.552
$ws.Range("I97").Value = 1112.3529
$ws.Range("J97").Value = 10303.417
$ws.Range("K97").Value = 1112.3529
$ws.Range("L97").Value = 10303.417
$ws.Range("M97").Value = -616.3529000000001
$ws.Range("N97").Value = -11295.417
# Row 130
$ws.Range("H130").Value = 44134.5
$ws.Range("J130").Value = 44134.5
$ws.Range("L130").Value = 44134.5
$ws.Range("N130").Value = -54174.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2314.7058
$ws.Range("I122").Value = 2288.889
$ws.Range("J122").Value = 2343.75
$ws.Range("K122").Value = 6866.667
$ws.Range("L122").Value = 7031.25
$ws.Range("M122").Value = -4416.667
$ws.Range("N122").Value = -11931.25
# Row 127
$ws.Range("H127").Value = 37357.5
$ws.Range("J127").Value = 37357.5
$ws.Range("L127").Value = 37357.5
$ws.Range("N127").Value = -47277.5
# Row 132
$ws.Range("H132").Value = 3590.721
$ws.Range("I132").Value = 3790.818
$ws.Range("J132").Value = 3381.0952
$ws.Range("K132").Value = 11372.454
$ws.Range("L132").Value = 10143.2856
$ws.Range("M132").Value = -8842.454000000002
$ws.Range("N132").Value = -15203.2856

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 46417.2
$ws.Range("J16").Value = 46417.2
$ws.Range("L16").Value = 46417.2
$ws.Range("N16").Value = -47001.2
# Row 120
$ws.Range("H120").Value = 42424.332
$ws.Range("J120").Value = 42424.332
$ws.Range("L120").Value = 42424.332
$ws.Range("N120").Value = -52100.332
# Row 121
$ws.Range("H121").Value = 43292
$ws.Range("J121").Value = 43292
$ws.Range("L121").Value = 43292
$ws.Range("N121").Value = -46786
# Row 126
$ws.Range("H126").Value = 500.36365
$ws.Range("I126").Value = 500.36365
$ws.Range("K126").Value = 1501.09095
$ws.Range("M126").Value = 968.90905
# Row 132
$ws.Range("H132").Value = 1052.3726
$ws.Range("I132").Value = 772.675
$ws.Range("K132").Value = 2318.025
$ws.Range("M132").Value = 211.9750000000004
